$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 461 (shifts all existing rows 461-551 down to 463-553).
$ws.Rows("461:462").Insert()

# Row 461: new price record for "Crespo record" variety.
$ws.Range("A461").Value = 10
$ws.Range("B461").Value = "Vega Modelo de Temuco"
$ws.Range("C461").Value = "La Araucanía"
$ws.Range("D461").Value = 44522
$ws.Range("E461").Value = 9
$ws.Range("F461").Value = 100112006
$ws.Range("G461").Value = "Repollo"
$ws.Range("H461").Value = "Crespo record"
$ws.Range("I461").Value = "Primera"
$ws.Range("J461").Value = 300
$ws.Range("K461").Value = 1000
$ws.Range("L461").Value = 1000
$ws.Range("M461").Value = 1000
$ws.Range("N461").Value = "$/unidad"
$ws.Range("O461").Value = "Región del Maule"
$ws.Range("P461").Value = 1000
$ws.Range("Q461").Value = 1
$ws.Range("R461").Value = "Hortaliza"

# Row 462: new price record for "Morada(o)" variety.
$ws.Range("A462").Value = 10
$ws.Range("B462").Value = "Vega Modelo de Temuco"
$ws.Range("C462").Value = "La Araucanía"
$ws.Range("D462").Value = 44522
$ws.Range("E462").Value = 9
$ws.Range("F462").Value = 100112006
$ws.Range("G462").Value = "Repollo"
$ws.Range("H462").Value = "Morada(o)"
$ws.Range("I462").Value = "Primera"
$ws.Range("J462").Value = 100
$ws.Range("K462").Value = 1000
$ws.Range("L462").Value = 1000
$ws.Range("M462").Value = 1000
$ws.Range("N462").Value = "$/unidad"
$ws.Range("O462").Value = "Región del Maule"
$ws.Range("P462").Value = 1000
$ws.Range("Q462").Value = 1
$ws.Range("R462").Value = "Hortaliza"
